# Update the quarterly income-statement database: shift the reporting
# window forward by one quarter (drop the oldest quarter column, add the
# newest quarter column) and refresh figures produced by the updated
# read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $row
        $ws.Range($addr).Value = $values[$i]
    }
}

# --- Row 8: quarter headers (دوره مالی) -------------------------------
Set-RowValues 8 @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

# --- Row 9: publish dates (تاریخ انتشار) -------------------------------
Set-RowValues 9 @(
    "1400-12-28 (3)",
    "1401-04-15 (10)",
    "1401-04-23 (3)",
    "1401-09-09 (5)",
    "1401-10-29 (3)",
    "1402-02-27 (9)",
    "1401-04-23",
    "1401-09-09 (3)",
    "1401-10-29",
    "1402-02-27 (2)"
)

# --- Row 11: فروش (Sales) ----------------------------------------------
Set-RowValues 11 @(50685, 98451, 75541, 93315, 70049, 112868, 85888, 74836, 87097, 72799)

# --- Row 12: بهای تمام شده کالای فروش رفته -----------------------------
Set-RowValues 12 @(-19940, -32650, -33515, -42606, -39763, -90084, -43474, -55736, -59036, -36051)

# --- Row 13: سود (زیان) ناخالص ------------------------------------------
Set-RowValues 13 @(30745, 65801, 42026, 50709, 30287, 22784, 42414, 19100, 28061, 36748)

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی --------------------------
Set-RowValues 14 @(-8125, -16770, -9723, -10823, -5142, -11619, -4001, -8299, -7725, -9416)

# Row 15 (هزینه کاهش ارزش دریافتنی ها) is unchanged - all dashes already.

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---------------------
Set-RowValues 16 @(7145, -5756, 502, 7485, 643, -2217, "-", 8294, 2036, 15544)

# --- Row 17: سود (زیان) عملیاتی -----------------------------------------
Set-RowValues 17 @(29765, 43276, 32805, 47371, 25788, 8949, 38413, 19094, 22371, 42876)

# Row 18 (هزینه های مالی) is unchanged - all dashes already.

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی --------------------
Set-RowValues 19 @(1869, 1829, 851, 2145, 1981, 2140, 1228, 3524, 833, 867)

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ----------
Set-RowValues 20 @(31634, 45105, 33655, 49516, 27769, 11089, 39640, 22618, 23204, 43743)

# --- Row 21: مالیات -------------------------------------------------------
Set-RowValues 21 @("-", "-", "-", "-", -202, -11563, -5257, -6693, -4760, -6292)

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم -------------------------
Set-RowValues 22 @(31634, 45105, 33655, 49516, 27567, -474, 34383, 15926, 18445, 37451)

# Row 23 (سود (زیان) عملیات متوقف شده پس از اثر مالیاتی) is unchanged - all dashes already.

# --- Row 24: سود (زیان) خالص ---------------------------------------------
Set-RowValues 24 @(31634, 45105, 33655, 49516, 27567, -474, 34383, 15926, 18445, 37451)

# --- Row 25: سود هر سهم پس از کسر مالیات ----------------------------------
Set-RowValues 25 @(0, 0, 0, 0, "-", 0, 0, 0, "-", 0)

# --- Row 26: سرمایه --------------------------------------------------------
Set-RowValues 26 @(22078, 24566, 25718, 22988, "-", 21784, 20419, 19345, "-", 13200)

# Row 27 (سود هر سهم بر اساس آخرین سرمایه) is unchanged - all zeros already.
